$wb = $excel.ActiveWorkbook

# "Hoja1" sheet - update the conversion text block
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.85 = 6854.04 pesos`n✅ 6854.04 pesos = 1.84 = 879.24 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# "tasas" sheet - update rate values
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 540.7
$wsTasas.Range("O10").Value = 3705.98
$wsTasas.Range("N12").Value = 3734
$wsTasas.Range("O12").Value = 479
